# Update "想去人数" (column F) and "最低票价" (column G) figures on the
# "展览" and "全部类型" sheets to match the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Row number -> hashtable of column letter -> new value.
# These updates apply to both the "展览" and "全部类型" sheets.
$commonUpdates = @{
    3  = @{ F = 413 }
    5  = @{ F = 124 }
    6  = @{ F = 38 }
    9  = @{ F = 119 }
    10 = @{ F = 1165 }
    11 = @{ F = 1476 }
    12 = @{ F = 327 }
    15 = @{ F = 114 }
    19 = @{ F = 259 }
    21 = @{ F = 305 }
    22 = @{ F = 1686 }
    25 = @{ F = 167 }
    26 = @{ F = 638 }
    28 = @{ F = 169 }
    29 = @{ F = 4013 }
    31 = @{ F = 472 }
    32 = @{ F = 243 }
    33 = @{ F = 1032 }
    34 = @{ F = 117 }
    36 = @{ F = 216; G = 60 }
    38 = @{ F = 138 }
}

# Extra, sheet-specific updates keyed by sheet name.
$extraUpdatesBySheet = @{
    "展览"     = @{}
    "全部类型" = @{ 7 = @{ F = 59 } }
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $commonUpdates.Keys) {
        $vals = $commonUpdates[$row]
        $ws.Cells.Item($row, 6).Value = $vals.F
        if ($vals.ContainsKey('G')) {
            $ws.Cells.Item($row, 7).Value = $vals.G
        }
    }

    $extraUpdates = $extraUpdatesBySheet[$sheetName]
    foreach ($row in $extraUpdates.Keys) {
        $vals = $extraUpdates[$row]
        $ws.Cells.Item($row, 6).Value = $vals.F
        if ($vals.ContainsKey('G')) {
            $ws.Cells.Item($row, 7).Value = $vals.G
        }
    }
}
